# Applies the re-ordering/update of species observation rows 29, 31, 34 and 35.
# The edit swaps the per-observation data between row 29 <-> row 35 and
# between row 31 <-> row 34 (only the columns that actually differ between
# the two rows of each pair are touched; columns that already hold identical
# values in both rows - dates, times, location names, reporter, etc. - are
# left untouched so no unrelated formatting/type coercion is introduced).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Pair 1 : row 29 <-> row 35
# ---------------------------------------------------------------------

# --- snapshot current ("before") values --------------------------------
$A29 = $ws.Range("A29").Value()
$B29 = $ws.Range("B29").Value()
$D29 = $ws.Range("D29").Value()
$E29 = $ws.Range("E29").Value()
$F29 = $ws.Range("F29").Value()
$G29 = $ws.Range("G29").Value()
$H29 = $ws.Range("H29").Value()
$I29 = $ws.Range("I29").Value()
$M29 = $ws.Range("M29").Value()
$Q29 = $ws.Range("Q29").Value()
$R29 = $ws.Range("R29").Value()

$A35 = $ws.Range("A35").Value()
$B35 = $ws.Range("B35").Value()
$D35 = $ws.Range("D35").Value()
$E35 = $ws.Range("E35").Value()
$F35 = $ws.Range("F35").Value()
$G35 = $ws.Range("G35").Value()
$H35 = $ws.Range("H35").Value()
$I35 = $ws.Range("I35").Value()
$J35 = $ws.Range("J35").Value()
$Q35 = $ws.Range("Q35").Value()
$R35 = $ws.Range("R35").Value()

# --- write row 29 with row 35's former values ---------------------------
$ws.Range("A29").Value = $A35
$ws.Range("B29").Value = $B35
$ws.Range("D29").Value = $D35
$ws.Range("E29").Value = $E35
$ws.Range("F29").Value = $F35
$ws.Range("G29").Value = $G35
$ws.Range("H29").Value = $H35
$ws.Range("I29").Value = $I35
$ws.Range("J29").Value = $J35
$ws.Range("M29").Value = ""
$ws.Range("Q29").Value = $Q35
$ws.Range("R29").Value = $R35
$ws.Range("AF29").Value = ""

# --- write row 35 with row 29's former values ---------------------------
$ws.Range("A35").Value = $A29
$ws.Range("B35").Value = $B29
$ws.Range("D35").Value = $D29
$ws.Range("E35").Value = $E29
$ws.Range("F35").Value = $F29
$ws.Range("G35").Value = $G29
$ws.Range("H35").Value = $H29
$ws.Range("I35").Value = $I29
$ws.Range("J35").Value = ""
$ws.Range("M35").Value = $M29
$ws.Range("Q35").Value = $Q29
$ws.Range("R35").Value = $R29
$ws.Range("AF35").Value = ""

# ---------------------------------------------------------------------
# Pair 2 : row 31 <-> row 34  (plain value swap, no cells added/removed)
# ---------------------------------------------------------------------

$A31 = $ws.Range("A31").Value()
$B31 = $ws.Range("B31").Value()
$D31 = $ws.Range("D31").Value()
$E31 = $ws.Range("E31").Value()
$F31 = $ws.Range("F31").Value()
$G31 = $ws.Range("G31").Value()
$H31 = $ws.Range("H31").Value()
$Q31 = $ws.Range("Q31").Value()
$R31 = $ws.Range("R31").Value()

$A34 = $ws.Range("A34").Value()
$B34 = $ws.Range("B34").Value()
$D34 = $ws.Range("D34").Value()
$E34 = $ws.Range("E34").Value()
$F34 = $ws.Range("F34").Value()
$G34 = $ws.Range("G34").Value()
$H34 = $ws.Range("H34").Value()
$Q34 = $ws.Range("Q34").Value()
$R34 = $ws.Range("R34").Value()

$ws.Range("A31").Value = $A34
$ws.Range("B31").Value = $B34
$ws.Range("D31").Value = $D34
$ws.Range("E31").Value = $E34
$ws.Range("F31").Value = $F34
$ws.Range("G31").Value = $G34
$ws.Range("H31").Value = $H34
$ws.Range("Q31").Value = $Q34
$ws.Range("R31").Value = $R34

$ws.Range("A34").Value = $A31
$ws.Range("B34").Value = $B31
$ws.Range("D34").Value = $D31
$ws.Range("E34").Value = $E31
$ws.Range("F34").Value = $F31
$ws.Range("G34").Value = $G31
$ws.Range("H34").Value = $H31
$ws.Range("Q34").Value = $Q31
$ws.Range("R34").Value = $R31
